$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A return target slider was added; this updates the "Opt Portfolio with
# View" column (D) with the resulting optimized weights.
$ws.Range("D2").Value = [double]"0.3467746958790734"
$ws.Range("D3").Value = [double]"0.3334176855470986"
$ws.Range("D4").Value = [double]"7.531898512759859E-19"
$ws.Range("D5").Value = [double]"2.188874988996542E-19"
$ws.Range("D6").Value = [double]"8.29629462328943E-19"
$ws.Range("D7").Value = [double]"0.3077857311009995"
$ws.Range("D8").Value = [double]"0.01202188747282855"
